# Insert a new data row at row 530 (pushing the existing rows 530..559
# down to 531..560), then populate the new row with its values.
#
# The new row is a near-duplicate of the row that used to be at 530
# (same Mercado/Región/Categoría/etc.) but with its own Fecha, price and
# $/Kg figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(530).Insert()

$ws.Cells.Item(530, 1).Value  = 5
$ws.Cells.Item(530, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(530, 3).Value  = "Maule"
$ws.Cells.Item(530, 4).Value  = 44931
$ws.Cells.Item(530, 5).Value  = 7
$ws.Cells.Item(530, 6).Value  = 100112043
$ws.Cells.Item(530, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(530, 8).Value  = "Sin especificar"
$ws.Cells.Item(530, 9).Value  = "Primera"
$ws.Cells.Item(530, 10).Value = 500
$ws.Cells.Item(530, 11).Value = 10000
$ws.Cells.Item(530, 12).Value = 10000
$ws.Cells.Item(530, 13).Value = 10000
$ws.Cells.Item(530, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(530, 15).Value = "Región del Maule"
$ws.Cells.Item(530, 16).Value = 125
$ws.Cells.Item(530, 17).Value = 80
$ws.Cells.Item(530, 18).Value = "Hortaliza"
